# Generate Report for Handoff
# Updates the localization-status workbook: refreshes the GUID / timestamps for the
# existing markdown entry and appends two new rows (for the two new .png assets)
# across the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "ae52835f-d048-4cc5-9112-25f87fdbc015"
$newGuid = "385e286e-d80d-4523-8cc5-101e6d8b4319"

$zhHash = "ddffd57b3b803030f67c005bceda1071856d0227"
$zhTime = "2016-03-18 07:27:32"
$deTime = "2016-03-18 07:27:34"
$overviewTime = "2016-27-18 07:27:34"

$png1 = "584007c9-9928-4934-bade-8c383208ec4e.png"
$png2 = "f1764c7f-2d50-48bc-ba2d-6db53af7fc6e.png"
$pngZh1 = "feb0c13cbbbb1827da643282c912022c8dff30ca.png"
$pngZh2 = "171cbf27d9c72264c8bdaf9dc300b6eec726498d.png"

$mdUrlBase   = "https://github.com/OpenLocalizationTest/oltest/blob/a6ff1fe4da77d0e348fd057a7bdbaab06c985dd3/e2e"
$zhUrlBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c440068960907256d586fa55c62d66c23a7c5353/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deUrlBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0938802f4431ffebe96eb23c0be688f8b69f406f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$mdFile   = "$newGuid.md"
$zhXlf    = "$newGuid.$zhHash.zh-cn.xlf"
$deXlf    = "$newGuid.$zhHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# -- Row 2: refresh GUID + hyperlink
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Range("A2").Value2 = $mdFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$mdUrlBase/$mdFile", "", "", $mdFile)
$wsOverview.Range("D2").Value2 = $overviewTime

# -- Row 3: new .png entry
$wsOverview.Range("B3").Value2 = "Ready for handoff"
$wsOverview.Range("C3").Value2 = "Ready for handoff"
$wsOverview.Range("D3").Value2 = $overviewTime
$wsOverview.Range("A3").Value2 = $png1
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$mdUrlBase/$png1", "", "", $png1)

# -- Row 4: new .png entry
$wsOverview.Range("B4").Value2 = "Ready for handoff"
$wsOverview.Range("C4").Value2 = "Ready for handoff"
$wsOverview.Range("D4").Value2 = $overviewTime
$wsOverview.Range("A4").Value2 = $png2
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$mdUrlBase/$png2", "", "", $png2)

Write-Host "Overview sheet updated"

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

# -- Row 2: refresh GUID, xlf hash and handoff timestamp
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("B2").Hyperlinks.Delete()
$wsZh.Range("D2").Hyperlinks.Delete()

$wsZh.Range("A2").Value2 = $mdFile
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$mdUrlBase/$mdFile", "", "", $mdFile)

$wsZh.Range("B2").Value2 = ".md"
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "$mdUrlBase/$mdFile", "", "", ".md")

$wsZh.Range("D2").Value2 = $zhXlf
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "$zhUrlBase/$zhXlf", "", "", $zhXlf)

$wsZh.Range("E2").Value2 = $zhTime
$wsZh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- Row 3: new .png entry (direct asset, Include)
$wsZh.Range("A3").Value2 = $png1
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$mdUrlBase/$png1", "", "", $png1)

$wsZh.Range("B3").Value2 = ".png"
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "$mdUrlBase/$png1", "", "", ".png")

$wsZh.Range("C3").Value2 = "Ready for handoff"

$wsZh.Range("D3").Value2 = $pngZh1
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "$zhUrlBase/$pngZh1", "", "", $pngZh1)

$wsZh.Range("E3").Value2 = $zhTime
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value2 = "IsDependency"
$wsZh.Range("J3").Value2 = "e2e\$mdFile"

# -- Row 4: new .png entry (dependency asset)
$wsZh.Range("A4").Value2 = $png2
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$mdUrlBase/$png2", "", "", $png2)

$wsZh.Range("B4").Value2 = ".png"
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), "$mdUrlBase/$png2", "", "", ".png")

$wsZh.Range("C4").Value2 = "Ready for handoff"

$wsZh.Range("D4").Value2 = $pngZh2
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "$zhUrlBase/$pngZh2", "", "", $pngZh2)

$wsZh.Range("E4").Value2 = $zhTime
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value2 = "IsDependency"
$wsZh.Range("J4").Value2 = "e2e\$mdFile"

Write-Host "zh-cn sheet updated"

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

# -- Row 2: refresh GUID, xlf hash and handoff timestamp
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("B2").Hyperlinks.Delete()
$wsDe.Range("D2").Hyperlinks.Delete()

$wsDe.Range("A2").Value2 = $mdFile
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$mdUrlBase/$mdFile", "", "", $mdFile)

$wsDe.Range("B2").Value2 = ".md"
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "$mdUrlBase/$mdFile", "", "", ".md")

$wsDe.Range("D2").Value2 = $deXlf
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "$deUrlBase/$deXlf", "", "", $deXlf)

$wsDe.Range("E2").Value2 = $deTime
$wsDe.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- Row 3: new .png entry (direct asset, Include)
$wsDe.Range("A3").Value2 = $png1
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$mdUrlBase/$png1", "", "", $png1)

$wsDe.Range("B3").Value2 = ".png"
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "$mdUrlBase/$png1", "", "", ".png")

$wsDe.Range("C3").Value2 = "Ready for handoff"

$wsDe.Range("D3").Value2 = $pngZh1
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "$deUrlBase/$pngZh1", "", "", $pngZh1)

$wsDe.Range("E3").Value2 = $deTime
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value2 = "IsDependency"
$wsDe.Range("J3").Value2 = "e2e\$mdFile"

# -- Row 4: new .png entry (dependency asset)
$wsDe.Range("A4").Value2 = $png2
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$mdUrlBase/$png2", "", "", $png2)

$wsDe.Range("B4").Value2 = ".png"
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), "$mdUrlBase/$png2", "", "", ".png")

$wsDe.Range("C4").Value2 = "Ready for handoff"

$wsDe.Range("D4").Value2 = $pngZh2
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "$deUrlBase/$pngZh2", "", "", $pngZh2)

$wsDe.Range("E4").Value2 = $deTime
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value2 = "IsDependency"
$wsDe.Range("J4").Value2 = "e2e\$mdFile"

Write-Host "de-de sheet updated"
